# Fruta / hortaliza, semanal
# Insert two new daily price rows (Early Glo / Super Queen, 10-kg trays from
# Provincia de Limarí, fecha 44524) ahead of the existing row 226 block,
# pushing the rest of the dataset down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 226-227; everything from the old row 226 onward
# shifts down to 228 onward (old row 303 becomes 305).
$ws.Range("A226:A227").EntireRow.Insert()

# Constant columns shared by every data row in this block.
$mercadoId = 10
$mercado = "Vega Modelo de Temuco"
$region = "La Araucanía"
$codreg = 9
$tipo = "Fruta"
$productoId = 100103
$producto = "Frutos de hueso (carozo)"
$categoriaId = 100103006
$categoria = "Nectarín"

# --- New row 226 ---
$ws.Cells.Item(226, 1).Value = $mercadoId
$ws.Cells.Item(226, 2).Value = $mercado
$ws.Cells.Item(226, 3).Value = $region
$ws.Cells.Item(226, 4).Value = 44524
$ws.Cells.Item(226, 5).Value = $codreg
$ws.Cells.Item(226, 6).Value = $tipo
$ws.Cells.Item(226, 7).Value = $productoId
$ws.Cells.Item(226, 8).Value = $producto
$ws.Cells.Item(226, 9).Value = $categoriaId
$ws.Cells.Item(226, 10).Value = $categoria
$ws.Cells.Item(226, 11).Value = "Early Glo"
$ws.Cells.Item(226, 12).Value = "Primera"
$ws.Cells.Item(226, 13).Value = 250
$ws.Cells.Item(226, 14).Value = 15000
$ws.Cells.Item(226, 15).Value = 15000
$ws.Cells.Item(226, 16).Value = 15000
$ws.Cells.Item(226, 17).Value = "$/bandeja 10 kilos granel"
$ws.Cells.Item(226, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(226, 19).Value = 1500
$ws.Cells.Item(226, 20).Value = 10

# --- New row 227 ---
$ws.Cells.Item(227, 1).Value = $mercadoId
$ws.Cells.Item(227, 2).Value = $mercado
$ws.Cells.Item(227, 3).Value = $region
$ws.Cells.Item(227, 4).Value = 44524
$ws.Cells.Item(227, 5).Value = $codreg
$ws.Cells.Item(227, 6).Value = $tipo
$ws.Cells.Item(227, 7).Value = $productoId
$ws.Cells.Item(227, 8).Value = $producto
$ws.Cells.Item(227, 9).Value = $categoriaId
$ws.Cells.Item(227, 10).Value = $categoria
$ws.Cells.Item(227, 11).Value = "Super Queen"
$ws.Cells.Item(227, 12).Value = "Primera"
$ws.Cells.Item(227, 13).Value = 250
$ws.Cells.Item(227, 14).Value = 20000
$ws.Cells.Item(227, 15).Value = 20000
$ws.Cells.Item(227, 16).Value = 20000
$ws.Cells.Item(227, 17).Value = "$/bandeja 10 kilos granel"
$ws.Cells.Item(227, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(227, 19).Value = 2000
$ws.Cells.Item(227, 20).Value = 10
